# Generate Report for handback
# Refresh the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# columns on the per-locale handback sheets.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
# Row 2: 2b2c6534-... entry
$zh.Range("D2").Value = "2016-02-15 08:17:43"
$zh.Range("G2").Value = "2016-02-15 08:18:34"
# Row 3: ad3a3400-... entry
$zh.Range("D3").Value = "2016-02-15 08:14:58"
$zh.Range("G3").Value = "2016-02-15 08:16:21"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
# Row 2: 2b2c6534-... entry
$de.Range("D2").Value = "2016-02-15 08:17:58"
$de.Range("G2").Value = "2016-02-15 08:19:00"
# Row 3: ad3a3400-... entry
$de.Range("D3").Value = "2016-02-15 08:15:13"
$de.Range("G3").Value = "2016-02-15 08:16:51"
